$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.070.50"
$ws.Range("E2").Value = "  +2.18%  "

$ws.Range("D3").Value = "1.905.54"
$ws.Range("E3").Value = "  +1.90%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.06"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.42%  "

$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4641"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -1.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4112"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +3.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.63"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07983"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -0.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.007"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +0.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.81"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -1.28%  "

$ws.Range("D13").Value = "1.910.16"
$ws.Range("E13").Value = "  +2.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.936"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -1.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.088"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -2.77%  "

$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.10"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -2.01%  "

$ws.Range("E18").Value = "  -0.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06579"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.51"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").Value = "29.130.73"
$ws.Range("E22").Value = "  +2.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.434"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +2.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.232"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -1.53%  "

$ws.Range("D26").Value = "2.139.78"
$ws.Range("E26").Value = "  +2.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.50"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -2.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.73"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.114"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.426"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -1.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.28"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -1.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9817"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +0.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09417"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -1.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.429"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +3.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.601"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +0.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.302"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -1.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06097"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -0.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02243"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.363"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("E40").Value = "  -0.40%  "

$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5799"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -2.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.17"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -1.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1824"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -2.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.266"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -1.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.325"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +11.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5506"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -1.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.03"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -1.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.915"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -2.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07051"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +1.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.81"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +19.16%  "
